# "Generate Report for Handoff"
# Status moves from "In Translation" -> "Ready for handoff" and the
# handoff/generate timestamps are refreshed for the zh-cn / de-de rows.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language Status cells + the generate date ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-19 16:43:48"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-19 16:43:44"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-19 16:43:48"

# The new "Ready for handoff" text is wider than "In Translation", so the
# Status columns re-autofit (Overview cols E/F, and col C on the language
# sheets) just like Excel does after a cell edit.
function Resize-ToFitHandoffStatus($col) {
  $col.ColumnWidth = 16.333333333333332
}

Resize-ToFitHandoffStatus $overview.Columns.Item(5)
Resize-ToFitHandoffStatus $overview.Columns.Item(6)
Resize-ToFitHandoffStatus $zhcn.Columns.Item(3)
Resize-ToFitHandoffStatus $dede.Columns.Item(3)
